# "remade external sorting algorithms logic"
#
# The sorting routine now places "Челябинск" ahead of "Уфа": the last two
# rows of the sorted-cities table swap their Word/Region data (Country stays
# "Russia" for both), i.e. row 11 becomes Челябинск/Europe and row 12
# becomes Уфа/Asia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Челябинск"
$ws.Range("B11").Value = "Europe"

$ws.Range("A12").Value = "Уфа"
$ws.Range("B12").Value = "Asia"
